$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute(
        $old, $true, $false, $false, $false, $false,
        $true, 1, $false, $new, 2
    ) | Out-Null
}

# --- Paragraph 0: add new trailing sentence + merge the "Al correr el programa..." runs ---
Replace-Exact "Al correr el programa lo primero que se despliega es un" "Al correr el programa lo primero que se despliega es un"
Replace-Exact "panel pidiendo ingresar un día." "panel pidiendo ingresar un día."
Replace-Exact "panel pidiendo ingresar un día." "panel pidiendo ingresar un día. Para simular un cambio de día, cierre el programa y ejecútelo de nuevo e inserte una nueva fecha en este panel."

# --- Paragraph 1: merge "frame de autenticac" and "ión en donde..." ---
Replace-Exact "frame de autenticac" "frame de autenticac"
Replace-Exact "ión en donde se puede iniciar sesión como administrador o como empleado." "ión en donde se puede iniciar sesión como administrador o como empleado."

# --- Paragraph 3: merge whole paragraph into one run ---
Replace-Exact "Para iniciar sesión como administrador ingrese las credenciales nombre: “root” y contraseña “Cookie”" "Para iniciar sesión como administrador ingrese las credenciales nombre: “root” y contraseña “Cookie”"

# --- Paragraph 5: merge "En caso de fallar la autenticac" and "ión, podrá..." ---
Replace-Exact "En caso de fallar la autenticac" "En caso de fallar la autenticac"
Replace-Exact "ión, podrá reintentarlo tantas veces como desee." "ión, podrá reintentarlo tantas veces como desee."

# --- Paragraph 8 ---
Replace-Exact " botón “Administrar usuarios”. Al hacer click en él se despliega una nueva ventana con un panel" " botón “Administrar usuarios”. Al hacer click en él se despliega una nueva ventana con un panel"
Replace-Exact " izquierdo para crear nuevos usuarios. Una vez se llenan todas las cajas del panel izquierdo se puede presionar el botón “Agregar usuario”. Al presionarlo se llena la lista del medio con el nuevo usuario, además de mostrar los antiguos usuarios existentes. Al hacer click en un usuario de la lista, su correspondiente información se muestra en el panel derecho. Una vez que un usuario ha sido seleccionado se puede presionar el botón “Quitar usuario”. Al confirmar la acción se elimina el usuario seleccionado de manera permanente. Para actualizar la información de un usuario existente se debe colocar el mismo nombre de usuario y cambiar los otros campos para finalmente darle “Agregar usuario” y actualizar su información en la lista." " izquierdo para crear nuevos usuarios. Una vez se llenan todas las cajas del panel izquierdo se puede presionar el botón “Agregar usuario”. Al presionarlo se llena la lista del medio con el nuevo usuario, además de mostrar los antiguos usuarios existentes. Al hacer click en un usuario de la lista, su correspondiente información se muestra en el panel derecho. Una vez que un usuario ha sido seleccionado se puede presionar el botón “Quitar usuario”. Al confirmar la acción se elimina el usuario seleccionado de manera permanente. Para actualizar la información de un usuario existente se debe colocar el mismo nombre de usuario y cambiar los otros campos para finalmente darle “Agregar usuario” y actualizar su información en la lista."

# --- Paragraph 10: merge whole paragraph ---
Replace-Exact "Como en todos los páneles, abajo a la izquierda hay un botón para regresar al menú principal (de Admin o de Empleado según corresponda)." "Como en todos los páneles, abajo a la izquierda hay un botón para regresar al menú principal (de Admin o de Empleado según corresponda)."

# --- Paragraph 12: merge whole paragraph ---
Replace-Exact "El botón “Tarifas” despliega un ventana para manejar las tarifas del hotel" "El botón “Tarifas” despliega un ventana para manejar las tarifas del hotel"

# --- Paragraph 14: merge runs 1..5 ---
Replace-Exact "Servicios” despliega una ventana que funciona muy similar al frame “Administrar Usuarios”. A la izquierda un panel para crear servicios. En el centro una lista con los servicios existentes. A la derecha un panel para consultar la información de un servicio seleccionado de la lista. Adicinalmente, un panel para cargar el servicio a una habitación en particular en caso de que se elija pagar luego, o donde también se puede elegir pagar en el momento." "Servicios” despliega una ventana que funciona muy similar al frame “Administrar Usuarios”. A la izquierda un panel para crear servicios. En el centro una lista con los servicios existentes. A la derecha un panel para consultar la información de un servicio seleccionado de la lista. Adicinalmente, un panel para cargar el servicio a una habitación en particular en caso de que se elija pagar luego, o donde también se puede elegir pagar en el momento."

# --- Paragraph 20 (Reservas) ---
Replace-Exact "espliega una " "espliega una "
Replace-Exact "entana que muestra" "entana que muestra"
Replace-Exact " la informac" " la informac"
Replace-Exact "ión de las reservas existentes. Arriba a la izquierda se puede buscar una reserva a partir de su número. Al ingresar el número y presionar “buscar reserva” se muestra en la tabla central la información de la reserva y en el panel superior la información del líder del grupo de la reserva." "ión de las reservas existentes. Arriba a la izquierda se puede buscar una reserva a partir de su número. Al ingresar el número y presionar “buscar reserva” se muestra en la tabla central la información de la reserva y en el panel superior la información del líder del grupo de la reserva."

# --- Paragraph 25: merge whole paragraph ---
Replace-Exact "Para añadir una habitación se debe presionar el botón “Agregar habitación”. Luego se debe elegir el tipo de habitación deseada. A continuación, se desplegarán todas las habitaciones que estén disponibles para las fechas seleccionadas y que sean del tipo pedido (ej. Estándar”), y, de entre la lista se debe hacer click en una de las habitaciones para que finalmente se cargue en la tabla central." "Para añadir una habitación se debe presionar el botón “Agregar habitación”. Luego se debe elegir el tipo de habitación deseada. A continuación, se desplegarán todas las habitaciones que estén disponibles para las fechas seleccionadas y que sean del tipo pedido (ej. Estándar”), y, de entre la lista se debe hacer click en una de las habitaciones para que finalmente se cargue en la tabla central."

# --- Paragraph 30: merge runs 0..2 ---
Replace-Exact "El botón “cerrar sesión” cierra la sesión del usuario actual y regresa al frame de autenticación" "El botón “cerrar sesión” cierra la sesión del usuario actual y regresa al frame de autenticación"

# --- Paragraph 32: merge whole paragraph ---
Replace-Exact "Por último, los botones “Check-in” y “Check-out” despliegan ventanas emergentes en dónde se puede colocar el número de un grupo para realizar la respectiva acción." "Por último, los botones “Check-in” y “Check-out” despliegan ventanas emergentes en dónde se puede colocar el número de un grupo para realizar la respectiva acción."
